$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Amira Sobhy, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Gehan Adel'
$ws.Range("G3").Value = 'Administrator, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Hend Mahmoud'
$ws.Range("G4").Value = 'Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Gehan Adel'
$ws.Range("G5").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G6").Value = 'Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef'
$ws.Range("G7").Value = 'Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Range("G11").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G13").Value = 'Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G19").Value = 'Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range("G20").Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range("G24").Value = 'Dr. Youstina Gamil, Dr. Sarah Mahdy'
$ws.Range("G30").Value = 'Dr. Yassmen Ahmad, Dr. Shorok Mohammad, Dr. Aya Hanafy, Dr. Wafaa Ebida'
